# Daily attendance processing - 2026-02-02 06:05:02
#
# In the "Recorded By" column (G), every cell whose value is the literal
# text "System, dnasr281@gmail.com" is updated so the two names are
# swapped: "dnasr281@gmail.com, System".
#
# We use Range.Find / Range.FindNext (rather than looping over every row
# with Cells.Item) so that only cells which already contain the target
# text are touched; this avoids materializing/disturbing the many blank
# "Recorded By" cells that exist elsewhere in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$searchText  = "System, dnasr281@gmail.com"
$replaceText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$colG = $ws.Range("G1:G" + $lastRow)

$replacedCount = 0
$firstAddress = $null

$found = $colG.Find($searchText)
while ($found -ne $null) {
    $addr = $found.Address()
    if ($firstAddress -eq $null) {
        $firstAddress = $addr
    } elseif ($addr -eq $firstAddress) {
        break
    }

    $found.Value = $replaceText
    $replacedCount = $replacedCount + 1

    $found = $colG.FindNext($found)
}

Write-Output "Updated $replacedCount 'Recorded By' cell(s) in column G."
